# General changes in infra and xls reading:
#  - add function to read addresses      -> Addresses sheet header D1: "adddress" -> "adddressLine"
#  - modify the function that read the data from sheet -> Registration sheet row heights
#  - add function to read products       -> rename "Products"/"Addresses" tabs to lowercase and
#                                            switch the active tab to "products"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Registration sheet: tweak a couple of row heights
# ---------------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Registration")
$wsReg.Rows.Item(2).RowHeight = 165
$wsReg.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------------
# 2) Addresses sheet: rename the "adddress" header to "adddressLine" and move
#    the active-cell selection onto that column
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("Addresses")
$wsAddr.Range("D1").Value = "adddressLine"
$wsAddr.Range("D2").Select()

# ---------------------------------------------------------------------------
# 3) Rename the "Products"/"Addresses" tabs to lower-case
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Products")
$wsProd.Name = "products"
$wsAddr.Name = "addresses"

# ---------------------------------------------------------------------------
# 4) Move the active tab from "cards" to "products"
# ---------------------------------------------------------------------------
$wsProd.Activate()
